$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 9, 10, 10, 10, 8, 5, 10, 6, 9, 10, 87),
    @(4, 9, 8, 7, 10, 8, 6, 9, 10, 8, 8, 83),
    @(8, 10, 10, 10, 10, 9, 10, 10, 10, 10, 10, 99),
    @(12, 6, 7, 7, 6, 7, 9, 8, 6, 9, 6, 81),
    @(13, 6, 7, 7, 6, 7, 9, 8, 6, 9, 6, 81),
    @(14, 9, 10, 8, 8, 7, 10, 10, 10, 10, 10, 92),
    @(16, 8, 10, 10, 8, 7, 10, 8, 8, 10, 8, 87)
)

$row = 2
foreach ($rowData in $data) {
    $col = 1
    foreach ($val in $rowData) {
        $ws.Cells.Item($row, $col).Value = $val
        $col++
    }
    $row++
}

$ws.Cells.Item(9, 11).Value = "Average"
$ws.Cells.Item(9, 12).Value = 87.14

$ws.Range("I16").Select()

$wb.Save()
